$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 181.8232256666666
$ws.Range("N2").Value = 545.4696769999999
$ws.Range("O2").Value = 0.5898296910336229
$ws.Range("P2").Value = 0.5898296910336229
$ws.Range("Q2").Value = 25983.87595455273
$ws.Range("R2").Value = 233854.8835909746
$ws.Range("S2").Value = 0.3289282093479995
$ws.Range("T2").Value = 0.3289282093479996

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.367786666666667
$ws.Range("N3").Value = 7.10336
$ws.Range("O3").Value = 0.007681036748263821
$ws.Range("P3").Value = 0.007681036748263821
$ws.Range("Q3").Value = 338.3741257912889
$ws.Range("R3").Value = 3045.3671321216
$ws.Range("S3").Value = 0.004283456228042914
$ws.Range("T3").Value = 0.004283456228042915

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 124.0729296666667
$ws.Range("N4").Value = 372.218789
$ws.Range("O4").Value = 0.4024892722181133
$ws.Range("P4").Value = 0.4024892722181133
$ws.Range("Q4").Value = 17730.93399897615
$ws.Range("R4").Value = 159578.4059907853
$ws.Range("S4").Value = 0.2244547495743762
$ws.Range("T4").Value = 0.2244547495743763

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 181.8232256666666
$ws.Range("N5").Value = 545.4696769999999
$ws.Range("O5").Value = 0.5898296910336229
$ws.Range("P5").Value = 0.5898296910336229
$ws.Range("Q5").Value = 11630.83337120923
$ws.Range("R5").Value = 104677.5003408831
$ws.Range("S5").Value = 0.1472339692780319
$ws.Range("T5").Value = 0.1472339692780319

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.367786666666667
$ws.Range("N6").Value = 7.10336
$ws.Range("O6").Value = 0.007681036748263821
$ws.Range("P6").Value = 0.007681036748263821
$ws.Range("Q6").Value = 151.46212524608
$ws.Range("R6").Value = 1363.15912721472
$ws.Range("S6").Value = 0.001917349271847427
$ws.Range("T6").Value = 0.001917349271847427

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 124.0729296666667
$ws.Range("N7").Value = 372.218789
$ws.Range("O7").Value = 0.4024892722181133
$ws.Range("P7").Value = 0.4024892722181133
$ws.Range("Q7").Value = 7936.673466987767
$ws.Range("R7").Value = 71430.0612028899
$ws.Range("S7").Value = 0.100469837380772
$ws.Range("T7").Value = 0.100469837380772

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 181.8232256666666
$ws.Range("N8").Value = 545.4696769999999
$ws.Range("O8").Value = 0.5898296910336229
$ws.Range("P8").Value = 0.5898296910336229
$ws.Range("Q8").Value = 8979.23151169035
$ws.Range("R8").Value = 80813.08360521316
$ws.Range("S8").Value = 0.1136675124075914
$ws.Range("T8").Value = 0.1136675124075914

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.367786666666667
$ws.Range("N9").Value = 7.10336
$ws.Range("O9").Value = 0.007681036748263821
$ws.Range("P9").Value = 0.007681036748263821
$ws.Range("Q9").Value = 116.9317317539556
$ws.Range("R9").Value = 1052.3855857856
$ws.Range("S9").Value = 0.001480231248373479
$ws.Range("T9").Value = 0.00148023124837348

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 124.0729296666667
$ws.Range("N10").Value = 372.218789
$ws.Range("O10").Value = 0.4024892722181133
$ws.Range("P10").Value = 0.4024892722181133
$ws.Range("Q10").Value = 6127.267601406966
$ws.Range("R10").Value = 55145.40841266269
$ws.Range("S10").Value = 0.07756468526296495
$ws.Range("T10").Value = 0.07756468526296495

Write-Output "done"
